$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "brand"
$ws.Range("J1").Value = "range"
$ws.Range("A2").Value = "asdasdasdsa"
$ws.Range("A8").Value = "hsfghfgh"
$ws.Range("A5").Value = "CM1008RWD-KALORIK"

$ws.Hyperlinks.Add($ws.Range("B5"), "")
$ws.Range("B5").WrapText = $true
$ws.Range("B5").VerticalAlignment = -4108

$a3 = $ws.Range("A3")
$a3.Font.Bold = $true
$a3.Font.Size = 24
$a3.VerticalAlignment = -4108

$a5 = $ws.Range("A5")
$a5.Font.ThemeColor = 1
$a5.VerticalAlignment = -4108
